# Auto-generated edit script applying value updates described in the commit diff.
# Each worksheet section updates the H-N ("price/profit") columns for specific leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder | Antidote
$ws.Range("H6").Value = 1992.75
$ws.Range("I6").Value = 185.5
$ws.Range("K6").Value = 556.5
$ws.Range("M6").Value = -444.5

# Row 17: One for the Road | Potion
$ws.Range("H17").Value = 1804.6666
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 2152.8572
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 6458.571599999999
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -6794.571599999999

# Row 39: Riches' Brew | Hi-Potion of Mind
$ws.Range("H39").Value = 584.6
$ws.Range("I39").Value = 224.33333
$ws.Range("J39").Value = 1125
$ws.Range("K39").Value = 672.99999
$ws.Range("L39").Value = 3375
$ws.Range("M39").Value = -376.99999
$ws.Range("N39").Value = -3967

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 2365.1133
$ws.Range("I98").Value = 1721.15
$ws.Range("K98").Value = 1721.15
$ws.Range("M98").Value = -223.1500000000001

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 1609.5454
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1628.3721
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4885.1163
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -7101.1163

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 2365.1133
$ws.Range("I122").Value = 1721.15
$ws.Range("K122").Value = 5163.450000000001
$ws.Range("M122").Value = -2713.450000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 36: Hot for Teacher | Heavy Iron Armor
$ws.Range("H36").Value = 23999.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 23999.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 23999.5
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -24691.5

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1506
$ws.Range("I61").Value = 1000.6774
$ws.Range("K61").Value = 1000.6774
$ws.Range("M61").Value = -788.6774

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 3343.282
$ws.Range("I74").Value = 3758.3704
$ws.Range("J74").Value = 2409.3333
$ws.Range("K74").Value = 3758.3704
$ws.Range("L74").Value = 2409.3333
$ws.Range("M74").Value = -2884.3704
$ws.Range("N74").Value = -4157.3333

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 3343.282
$ws.Range("I77").Value = 3758.3704
$ws.Range("J77").Value = 2409.3333
$ws.Range("K77").Value = 18791.852
$ws.Range("L77").Value = 12046.6665
$ws.Range("M77").Value = -14423.852
$ws.Range("N77").Value = -20782.6665

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 1833.3334
$ws.Range("I102").Value = 1466.6666
$ws.Range("J102").Value = 2200
$ws.Range("K102").Value = 1466.6666
$ws.Range("L102").Value = 2200
$ws.Range("M102").Value = 155.3334
$ws.Range("N102").Value = -5444

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2403.8948
$ws.Range("I132").Value = 1743.9149
$ws.Range("J132").Value = 5505.8
$ws.Range("K132").Value = 5231.7447
$ws.Range("L132").Value = 16517.4
$ws.Range("M132").Value = -2701.7447
$ws.Range("N132").Value = -21577.4

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1506
$ws.Range("I136").Value = 1000.6774
$ws.Range("K136").Value = 3002.0322
$ws.Range("M136").Value = -452.0322000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 4971.054
$ws.Range("I20").Value = 1091.409
$ws.Range("J20").Value = 10661.2
$ws.Range("K20").Value = 1091.409
$ws.Range("L20").Value = 10661.2
$ws.Range("M20").Value = -844.4090000000001
$ws.Range("N20").Value = -11155.2

# Row 59: Pop That Top | Cobalt Raising Hammer
$ws.Range("H59").Value = 85950
$ws.Range("J59").Value = 85950
$ws.Range("L59").Value = 85950
$ws.Range("N59").Value = -87644

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 1644.4375
$ws.Range("I94").Value = 1721.5
$ws.Range("J94").Value = 1105
$ws.Range("K94").Value = 1721.5
$ws.Range("L94").Value = 1105
$ws.Range("M94").Value = -1270.5
$ws.Range("N94").Value = -2007

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2830.5938
$ws.Range("J105").Value = 2729.6667
$ws.Range("L105").Value = 2729.6667
$ws.Range("N105").Value = -6223.6667

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2226.7612
$ws.Range("I134").Value = 1196.4736
$ws.Range("J134").Value = 8099.4
$ws.Range("K134").Value = 3589.4208
$ws.Range("L134").Value = 24298.2
$ws.Range("M134").Value = -1054.4208
$ws.Range("N134").Value = -29368.2

# Row 137: Dagger Swagger | Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 32959.43
$ws.Range("J137").Value = 32959.43
$ws.Range("L137").Value = 32959.43
$ws.Range("N137").Value = -43159.43

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 5900.8
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5900.8
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5900.8
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6490.8

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 5900.8
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5900.8
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5900.8
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6304.8

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1848.2667
$ws.Range("I58").Value = 1510.8823
$ws.Range("J58").Value = 5125.7144
$ws.Range("K58").Value = 1510.8823
$ws.Range("L58").Value = 5125.7144
$ws.Range("M58").Value = -1307.8823
$ws.Range("N58").Value = -5531.7144

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 3771.5483
$ws.Range("I99").Value = 3011.2222
$ws.Range("K99").Value = 3011.2222
$ws.Range("M99").Value = -1513.2222

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 3771.5483
$ws.Range("I126").Value = 3011.2222
$ws.Range("K126").Value = 9033.6666
$ws.Range("M126").Value = -6563.6666

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1848.2667
$ws.Range("I136").Value = 1510.8823
$ws.Range("J136").Value = 5125.7144
$ws.Range("K136").Value = 4532.6469
$ws.Range("L136").Value = 15377.1432
$ws.Range("M136").Value = -1982.6469
$ws.Range("N136").Value = -20477.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa | Grilled Trout
$ws.Range("H3").Value = 4369.909
$ws.Range("I3").Value = 2484.6428
$ws.Range("J3").Value = 7669.125
$ws.Range("K3").Value = 7453.928400000001
$ws.Range("L3").Value = 23007.375
$ws.Range("M3").Value = -7341.928400000001
$ws.Range("N3").Value = -23231.375

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 657.86365
$ws.Range("I113").Value = 669.875
$ws.Range("J113").Value = 625.8333
$ws.Range("K113").Value = 2009.625
$ws.Range("L113").Value = 1877.4999
$ws.Range("M113").Value = 160.375
$ws.Range("N113").Value = -6217.4999

# Row 115: Mixology | Blood Tomato Juice
$ws.Range("H115").Value = 51260.5
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 51260.5
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 153781.5
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -156131.5

# Row 119: Super Dark Times | Risotto al Nero
$ws.Range("H119").Value = 4216.8
$ws.Range("I119").Value = 368
$ws.Range("J119").Value = 9990
$ws.Range("K119").Value = 1104
$ws.Range("L119").Value = 29970
$ws.Range("M119").Value = 3734
$ws.Range("N119").Value = -39646

$ws = $wb.Worksheets.Item("GSM")
# Row 46: Burning the Midnight Oil | Fire Brand
$ws.Range("H46").Value = 34073.6
$ws.Range("J46").Value = 34073.6
$ws.Range("L46").Value = 34073.6
$ws.Range("N46").Value = -34385.6

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 950.03705
$ws.Range("I97").Value = 919.5714
$ws.Range("J97").Value = 1056.6666
$ws.Range("K97").Value = 919.5714
$ws.Range("L97").Value = 1056.6666
$ws.Range("M97").Value = -423.5714
$ws.Range("N97").Value = -2048.6666

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 2436.075
$ws.Range("I102").Value = 1897.4642
$ws.Range("J102").Value = 3692.8333
$ws.Range("K102").Value = 1897.4642
$ws.Range("L102").Value = 3692.8333
$ws.Range("M102").Value = -275.4641999999999
$ws.Range("N102").Value = -6936.8333

# Row 109: You're My Wonderhall | Hematite Earrings of Healing
$ws.Range("H109").Value = 30285
$ws.Range("J109").Value = 30285
$ws.Range("L109").Value = 30285
$ws.Range("N109").Value = -32365

# Row 136: Shiny and Good | Pink Beryl
$ws.Range("H136").Value = 11229.3125
$ws.Range("I136").Value = 590
$ws.Range("J136").Value = 11572.517
$ws.Range("K136").Value = 1770
$ws.Range("L136").Value = 34717.551
$ws.Range("N136").Value = -39817.551
$ws.Range("M136").Value = 780

# Row 137: Sew Excited | Cobalt Tungsten Needle
$ws.Range("H137").Value = 52700
$ws.Range("J137").Value = 52700
$ws.Range("L137").Value = 52700
$ws.Range("N137").Value = -62900

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 3500.7222
$ws.Range("I40").Value = 3311
$ws.Range("J40").Value = 4286.7144
$ws.Range("K40").Value = 3311
$ws.Range("L40").Value = 4286.7144
$ws.Range("M40").Value = -3175
$ws.Range("N40").Value = -4558.7144

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 791.2105
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 791.2105
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Clothing the Naked Truth | Ramie Robe of Casting
$ws.Range("H74").Value = 7813.5
$ws.Range("J74").Value = 7813.5
$ws.Range("L74").Value = 7813.5
$ws.Range("N74").Value = -9685.5

# Row 77: When in Robes (L) | Ramie Robe of Casting
$ws.Range("H77").Value = 7813.5
$ws.Range("J77").Value = 7813.5
$ws.Range("L77").Value = 23440.5
$ws.Range("N77").Value = -32800.5

# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 7937680
$ws.Range("I81").Value = 8929702
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 17859404
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -17858343
$ws.Range("N81").Value = -5122

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 7937680
$ws.Range("I84").Value = 8929702
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 89297020
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -89291716
$ws.Range("N84").Value = -25608

# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 325.86667
$ws.Range("I113").Value = 277.7143
$ws.Range("K113").Value = 833.1428999999999
$ws.Range("M113").Value = 1336.8571

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2247.18
$ws.Range("I122").Value = 1567.9722
$ws.Range("J122").Value = 3993.7144
$ws.Range("K122").Value = 4703.9166
$ws.Range("L122").Value = 11981.1432
$ws.Range("M122").Value = -2253.9166
$ws.Range("N122").Value = -16881.1432
